$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Diode change: update placement coordinates for the affected parts ---
# C3 (row 3) — Mid X / Mid Y
$ws.Range("B3").Value = 65.45
$ws.Range("C3").Value = -52.9648

# D1 (row 7) — Mid X / Mid Y
$ws.Range("B7").Value = 65.775
$ws.Range("C7").Value = -54.225

# R2 (row 16) — Mid X / Mid Y
$ws.Range("B16").Value = 65.4602
$ws.Range("C16").Value = -51.898

# --- Materialize blank placeholder cells G3:H3 (extends used range to col H) ---
$touch = $ws.Range("G3:H3")
$touch.Font.Name = "宋体"
$touch.Font.Family = 0

# --- Restore the active selection to where the editor left off ---
$ws.Range("D31").Select()
